# Apply the split-plot-a fix (work-in-progress) edit described in the commit:
# "Starting fix on split plot analysis. Error found by Zhimin Pan"
#
# The edit touches the three data sheets (TP, FP, TRUTH) and relabels / retypes
# a number of cells. Many of the resulting values are stored as *text* rather
# than numbers (this mirrors the underlying OOXML diff, where those cells
# switch from numeric <v> cells to shared-string <v t="s"> cells).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "TP"
# ---------------------------------------------------------------------------
$tp = $wb.Worksheets.Item("TP")

# Header E1: "TP_Rating" -> "LL_Rating"
$tp.Range("E1").Value = "LL_Rating"

# Column B (rows 2-11): numeric 4 -> text "1"
$tp.Range("B2").Value = "1"
$tp.Range("B3").Value = "1"
$tp.Range("B4").Value = "1"
$tp.Range("B5").Value = "1"
$tp.Range("B6").Value = "1"
$tp.Range("B7").Value = "1"
$tp.Range("B8").Value = "1"
$tp.Range("B9").Value = "1"
$tp.Range("B10").Value = "1"
$tp.Range("B11").Value = "1"

# Column E: a few rating values become text labels
$tp.Range("E2").Value = "6"
$tp.Range("E3").Value = "2"
$tp.Range("E12").Value = "1"
$tp.Range("E13").Value = "2"

# ---------------------------------------------------------------------------
# Sheet 2: "FP"
# ---------------------------------------------------------------------------
$fp = $wb.Worksheets.Item("FP")

# Header D1: "FP_Rating" -> "NL_Rating"
$fp.Range("D1").Value = "NL_Rating"

# Column B (rows 2-11): numeric 4 -> text "1"
$fp.Range("B2").Value = "1"
$fp.Range("B3").Value = "1"
$fp.Range("B4").Value = "1"
$fp.Range("B5").Value = "1"
$fp.Range("B6").Value = "1"
$fp.Range("B7").Value = "1"
$fp.Range("B8").Value = "1"
$fp.Range("B9").Value = "1"
$fp.Range("B10").Value = "1"
$fp.Range("B11").Value = "1"

# ---------------------------------------------------------------------------
# Sheet 3: "TRUTH"
# ---------------------------------------------------------------------------
$truth = $wb.Worksheets.Item("TRUTH")

# Column E (rows 2-11): "(4)" -> "(1)"
$truth.Range("E2").Value = "(1)"
$truth.Range("E3").Value = "(1)"
$truth.Range("E4").Value = "(1)"
$truth.Range("E5").Value = "(1)"
$truth.Range("E6").Value = "(1)"
$truth.Range("E7").Value = "(1)"
$truth.Range("E8").Value = "(1)"
$truth.Range("E9").Value = "(1)"
$truth.Range("E10").Value = "(1)"
$truth.Range("E11").Value = "(1)"

# F3: "split-plot-a" -> "SPLIT-PLOT-A"
$truth.Range("F3").Value = "SPLIT-PLOT-A"
